$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" date column (C) for rows 2-7 from 2023-09-10 to 2023-09-11
foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = 45180
}
